# Updates cryptos price (D) / volume-change (E) columns to the latest
# scrape, mirroring the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.029.99"
$ws.Cells.Item(2, 5).Value = "  -0.42%  "

$ws.Cells.Item(3, 4).Value = "1.645.06"
$ws.Cells.Item(3, 5).Value = "  -1.48%  "

$ws.Cells.Item(4, 5).Value = "  -0.19%  "

$ws.Cells.Item(5, 4).Value = "'215.27"
$ws.Cells.Item(5, 5).Value = "  +2.31%  "

$ws.Cells.Item(6, 4).Value = "'0.5223"
$ws.Cells.Item(6, 5).Value = "  +0.11%  "

$ws.Cells.Item(7, 4).Value = "'1.001"
$ws.Cells.Item(7, 5).Value = "  -0.20%  "

$ws.Cells.Item(8, 4).Value = "'0.2609"
$ws.Cells.Item(8, 5).Value = "  -0.18%  "

$ws.Cells.Item(9, 4).Value = "'0.06354"
$ws.Cells.Item(9, 5).Value = "  +0.37%  "

$ws.Cells.Item(10, 4).Value = "'20.80"
$ws.Cells.Item(10, 5).Value = "  -1.67%  "

$ws.Cells.Item(11, 4).Value = "'0.07657"
$ws.Cells.Item(11, 5).Value = "  +1.45%  "

$ws.Cells.Item(12, 4).Value = "1.643.59"
$ws.Cells.Item(12, 5).Value = "  -1.75%  "

$ws.Cells.Item(13, 4).Value = "'4.419"
$ws.Cells.Item(13, 5).Value = "  -0.09%  "

$ws.Cells.Item(14, 4).Value = "1.867.53"
$ws.Cells.Item(14, 5).Value = "  -1.63%  "

$ws.Cells.Item(15, 4).Value = "'0.5536"
$ws.Cells.Item(15, 5).Value = "  +1.68%  "

$ws.Cells.Item(16, 4).Value = "0.0₅8295"
$ws.Cells.Item(16, 5).Value = "  +3.22%  "

$ws.Cells.Item(17, 4).Value = "'64.84"
$ws.Cells.Item(17, 5).Value = "  -2.45%  "

$ws.Cells.Item(18, 4).Value = "26.031.79"
$ws.Cells.Item(18, 5).Value = "  -0.48%  "

$ws.Cells.Item(19, 4).Value = "'1.001"
$ws.Cells.Item(19, 5).Value = "  -0.19%  "

$ws.Cells.Item(20, 4).Value = "'4.712"
$ws.Cells.Item(20, 5).Value = "  -0.70%  "

$ws.Cells.Item(21, 4).Value = "'188.18"
$ws.Cells.Item(21, 5).Value = "  +0.47%  "

$ws.Cells.Item(22, 4).Value = "'10.18"
$ws.Cells.Item(22, 5).Value = "  -1.05%  "

$ws.Cells.Item(23, 4).Value = "'6.253"
$ws.Cells.Item(23, 5).Value = "  +0.12%  "

$ws.Cells.Item(24, 5).Value = "  -0.21%  "

$ws.Cells.Item(25, 4).Value = "'145.32"
$ws.Cells.Item(25, 5).Value = "  -2.79%  "

$ws.Cells.Item(26, 4).Value = "'0.1217"
$ws.Cells.Item(26, 5).Value = "  -1.70%  "

$ws.Cells.Item(27, 4).Value = "'7.400"
$ws.Cells.Item(27, 5).Value = "  -1.11%  "

$ws.Cells.Item(28, 4).Value = "'15.82"

$ws.Cells.Item(29, 4).Value = "'1.390"
$ws.Cells.Item(29, 5).Value = "  +1.32%  "

$ws.Cells.Item(30, 4).Value = "'0.05950"
$ws.Cells.Item(30, 5).Value = "  -5.67%  "

$ws.Cells.Item(31, 5).Value = "  -1.27%  "

$ws.Cells.Item(32, 4).Value = "'3.406"
$ws.Cells.Item(32, 5).Value = "  -2.92%  "

$ws.Cells.Item(33, 4).Value = "'3.396"
$ws.Cells.Item(33, 5).Value = "  -0.82%  "

$ws.Cells.Item(34, 4).Value = "'1.651"
$ws.Cells.Item(34, 5).Value = "  +0.16%  "

$ws.Cells.Item(35, 4).Value = "'0.9939"
$ws.Cells.Item(35, 5).Value = "  -0.95%  "

$ws.Cells.Item(36, 4).Value = "'2.393"
$ws.Cells.Item(36, 5).Value = "  -0.22%  "

$ws.Cells.Item(37, 4).Value = "'2.750"
$ws.Cells.Item(37, 5).Value = "  -0.45%  "

$ws.Cells.Item(38, 4).Value = "'0.5625"
$ws.Cells.Item(38, 5).Value = "  -6.37%  "

$ws.Cells.Item(39, 4).Value = "'0.01608"
$ws.Cells.Item(39, 5).Value = "  -0.30%  "

$ws.Cells.Item(40, 4).Value = "'5.852"
$ws.Cells.Item(40, 5).Value = "  -3.45%  "

$ws.Cells.Item(41, 4).Value = "'0.8547"
$ws.Cells.Item(41, 5).Value = "  -1.20%  "

$ws.Cells.Item(42, 5).Value = "  -0.28%  "

$ws.Cells.Item(43, 4).Value = "1.024.99"
$ws.Cells.Item(43, 5).Value = "  -8.18%  "

$ws.Cells.Item(44, 4).Value = "'98.50"
$ws.Cells.Item(44, 5).Value = "  -2.18%  "

$ws.Cells.Item(45, 4).Value = "1.794.84"
$ws.Cells.Item(45, 5).Value = "  -1.47%  "

$ws.Cells.Item(46, 5).Value = "  +0.76%  "

$ws.Cells.Item(47, 4).Value = "'55.68"
$ws.Cells.Item(47, 5).Value = "  +0.25%  "

$ws.Cells.Item(48, 4).Value = "'1.006"
$ws.Cells.Item(48, 5).Value = "  +0.55%  "

$ws.Cells.Item(49, 4).Value = "'8.070"
$ws.Cells.Item(49, 5).Value = "  -0.04%  "

$ws.Cells.Item(50, 4).Value = "'0.05146"
$ws.Cells.Item(50, 5).Value = "  -2.04%  "

$ws.Cells.Item(51, 4).Value = "'0.4216"
$ws.Cells.Item(51, 5).Value = "  -0.58%  "
